# Insert a new data row at row 588 (pushing the existing rows 588-629 down
# to 589-630) and populate it with the new record:
#   2026/01/09  金  20  24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(588).Insert()

# Column A holds a date formatted as plain text (e.g. "2026/01/09"), not a
# real Excel date. Assigning that text straight to Range.Value would make
# Excel auto-convert it into a date serial number (and stamp a date
# NumberFormat/style onto the cell), which would not match the source
# data. Instead, compute the literal string via a formula in a scratch
# cell and paste back only the resulting value - this keeps the cell as
# plain text with the default (unstyled) formatting, exactly like its
# neighboring cells.
$ws.Range("Z1").Formula = "=""2026/01/09"""
$ws.Range("Z1").Copy()
$ws.Range("A588").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("B588").Value = "金"
$ws.Range("C588").Value = 20
$ws.Range("D588").Value = 24
